$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column M data mirroring the existing pattern (column L) one column to the right
$ws.Range("L3").Copy() | Out-Null
$ws.Range("M3").PasteSpecial(-4122) | Out-Null

$ws.Range("L4").Copy() | Out-Null
$ws.Range("M4").PasteSpecial(-4122) | Out-Null
$ws.Range("M4").Value = 2022

$ws.Range("L5").Copy() | Out-Null
$ws.Range("M5").PasteSpecial(-4122) | Out-Null
$ws.Range("M5").Value = 373

$excel.CutCopyMode = 0

# Update the selected cell/range to match the new state
$ws.Range("O4").Select() | Out-Null
